# Update cached calculated values in the FFXIV leve-profit sheets
# (market price / profit figures refreshed by the scheduled data runner).
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 315256.66
$ws.Range("J17").Value = 325372.6
$ws.Range("L17").Value = 976117.7999999999
$ws.Range("N17").Value = -976453.7999999999
$ws.Range("H19").Value = 1195.3125
$ws.Range("I19").Value = 1460.4445
$ws.Range("J19").Value = 854.4286
$ws.Range("K19").Value = 1460.4445
$ws.Range("L19").Value = 854.4286
$ws.Range("M19").Value = -1285.4445
$ws.Range("N19").Value = -1204.4286
$ws.Range("H33").Value = 509.2857
$ws.Range("I33").Value = 499.22223
$ws.Range("J33").Value = 527.4
$ws.Range("K33").Value = 499.22223
$ws.Range("L33").Value = 527.4
$ws.Range("M33").Value = -270.22223
$ws.Range("N33").Value = -985.4
$ws.Range("H64").Value = 6427.5713
$ws.Range("I64").Value = 5999
$ws.Range("K64").Value = 5999
$ws.Range("M64").Value = -5751
$ws.Range("H67").Value = 6427.5713
$ws.Range("I67").Value = 5999
$ws.Range("K67").Value = 5999
$ws.Range("M67").Value = -5141
$ws.Range("H86").Value = 3987.4614
$ws.Range("J86").Value = 4607.5
$ws.Range("L86").Value = 4607.5
$ws.Range("N86").Value = -6853.5
$ws.Range("H89").Value = 3987.4614
$ws.Range("J89").Value = 4607.5
$ws.Range("L89").Value = 23037.5
$ws.Range("N89").Value = -34269.5
$ws.Range("H111").Value = 16619.75
$ws.Range("J111").Value = 6828
$ws.Range("L111").Value = 20484
$ws.Range("N111").Value = -26618
$ws.Range("H121").Value = 2165
$ws.Range("H132").Value = 18524742
$ws.Range("I132").Value = 21279830
$ws.Range("J132").Value = 26289.715
$ws.Range("K132").Value = 63839490
$ws.Range("L132").Value = 78869.145
$ws.Range("M132").Value = -63836960
$ws.Range("N132").Value = -83929.145
$ws.Range("H138").Value = 1002351.8
$ws.Range("J138").Value = 1670399.9
$ws.Range("L138").Value = 5011199.699999999
$ws.Range("N138").Value = -5021479.699999999

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2089.0588
$ws.Range("I61").Value = 1904.5
$ws.Range("J61").Value = 3473.25
$ws.Range("K61").Value = 1904.5
$ws.Range("L61").Value = 3473.25
$ws.Range("M61").Value = -1692.5
$ws.Range("N61").Value = -3897.25
$ws.Range("H74").Value = 2393.2693
$ws.Range("I74").Value = 2005.5
$ws.Range("K74").Value = 2005.5
$ws.Range("M74").Value = -1131.5
$ws.Range("H77").Value = 2393.2693
$ws.Range("I77").Value = 2005.5
$ws.Range("K77").Value = 10027.5
$ws.Range("M77").Value = -5659.5
$ws.Range("H122").Value = 4548.5835
$ws.Range("I122").Value = 4508.5
$ws.Range("K122").Value = 13525.5
$ws.Range("M122").Value = -11075.5
$ws.Range("H132").Value = 6624.2856
$ws.Range("I132").Value = 2521.9092
$ws.Range("K132").Value = 7565.7276
$ws.Range("M132").Value = -5035.7276
$ws.Range("H136").Value = 2089.0588
$ws.Range("I136").Value = 1904.5
$ws.Range("J136").Value = 3473.25
$ws.Range("K136").Value = 5713.5
$ws.Range("L136").Value = 10419.75
$ws.Range("M136").Value = -3163.5
$ws.Range("N136").Value = -15519.75

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3343.7083
$ws.Range("I105").Value = 2445.9473
$ws.Range("K105").Value = 2445.9473
$ws.Range("M105").Value = -698.9472999999998
$ws.Range("H107").Value = 804.6896400000001
$ws.Range("I107").Value = 644.92
$ws.Range("K107").Value = 644.92
$ws.Range("M107").Value = 1275.08
$ws.Range("H134").Value = 4521.351
$ws.Range("I134").Value = 1710.8572
$ws.Range("K134").Value = 5132.571599999999
$ws.Range("M134").Value = -2597.571599999999

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2527.15
$ws.Range("I58").Value = 2432.75
$ws.Range("K58").Value = 2432.75
$ws.Range("M58").Value = -2229.75
$ws.Range("H99").Value = 1164.8
$ws.Range("I99").Value = 1164.8
$ws.Range("K99").Value = 1164.8
$ws.Range("M99").Value = 333.2
$ws.Range("H107").Value = 3029.182
$ws.Range("I107").Value = 2512.875
$ws.Range("K107").Value = 2512.875
$ws.Range("M107").Value = -592.875
$ws.Range("H126").Value = 1164.8
$ws.Range("I126").Value = 1164.8
$ws.Range("K126").Value = 3494.4
$ws.Range("M126").Value = -1024.4
$ws.Range("H134").Value = 3005.7144
$ws.Range("I134").Value = 2116.25
$ws.Range("K134").Value = 6348.75
$ws.Range("M134").Value = -3813.75
$ws.Range("H136").Value = 2527.15
$ws.Range("I136").Value = 2432.75
$ws.Range("K136").Value = 7298.25
$ws.Range("M136").Value = -4748.25

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 90.93939
$ws.Range("I2").Value = 83.541664
$ws.Range("J2").Value = 110.666664
$ws.Range("K2").Value = 501.249984
$ws.Range("L2").Value = 663.999984
$ws.Range("M2").Value = -388.249984
$ws.Range("N2").Value = -889.999984
$ws.Range("H86").Value = 650.36365
$ws.Range("I86").Value = 650.44446
$ws.Range("J86").Value = 650
$ws.Range("K86").Value = 1951.33338
$ws.Range("L86").Value = 1950
$ws.Range("M86").Value = -765.33338
$ws.Range("N86").Value = -4322
$ws.Range("H89").Value = 650.36365
$ws.Range("I89").Value = 650.44446
$ws.Range("J89").Value = 650
$ws.Range("K89").Value = 5854.00014
$ws.Range("L89").Value = 5850
$ws.Range("M89").Value = 73.9998599999999
$ws.Range("N89").Value = -17706
$ws.Range("H97").Value = 1572.7778
$ws.Range("I97").Value = 1909.1666
$ws.Range("J97").Value = 900
$ws.Range("K97").Value = 5727.4998
$ws.Range("L97").Value = 2700
$ws.Range("M97").Value = -5231.4998
$ws.Range("N97").Value = -3692
$ws.Range("H128").Value = 96999.336
$ws.Range("I128").Value = 96999.336
$ws.Range("K128").Value = 290998.008
$ws.Range("M128").Value = -286018.008

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2308.4285
$ws.Range("J97").Value = 2252.5
$ws.Range("L97").Value = 2252.5
$ws.Range("N97").Value = -3244.5
$ws.Range("H122").Value = 1996
$ws.Range("I122").Value = 1829.05
$ws.Range("K122").Value = 5487.15
$ws.Range("M122").Value = -3037.15
$ws.Range("H132").Value = 5218.6
$ws.Range("I132").Value = 5218.6
$ws.Range("K132").Value = 15655.8
$ws.Range("M132").Value = -13125.8

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5322.8184
$ws.Range("I7").Value = 3359
$ws.Range("K7").Value = 3359
$ws.Range("M7").Value = -3247
$ws.Range("H16").Value = 1302.8462
$ws.Range("I16").Value = 1302.8462
$ws.Range("K16").Value = 1302.8462
$ws.Range("M16").Value = -1132.8462
$ws.Range("H40").Value = 3585
$ws.Range("I40").Value = 3180.7
$ws.Range("K40").Value = 3180.7
$ws.Range("M40").Value = -3044.7
$ws.Range("H126").Value = 5322.8184
$ws.Range("I126").Value = 3359
$ws.Range("K126").Value = 10077
$ws.Range("M126").Value = -7607

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 34989
$ws.Range("J47").Value = 34989
$ws.Range("L47").Value = 34989
$ws.Range("N47").Value = -36133
$ws.Range("H81").Value = 56690.95
$ws.Range("I81").Value = 147171.58
$ws.Range("J81").Value = 7970.615
$ws.Range("K81").Value = 294343.16
$ws.Range("L81").Value = 15941.23
$ws.Range("M81").Value = -293282.16
$ws.Range("N81").Value = -18063.23
$ws.Range("H84").Value = 56690.95
$ws.Range("I84").Value = 147171.58
$ws.Range("J84").Value = 7970.615
$ws.Range("K84").Value = 1471715.8
$ws.Range("L84").Value = 79706.14999999999
$ws.Range("M84").Value = -1466411.8
$ws.Range("N84").Value = -90314.14999999999
$ws.Range("H100").Value = 1129.2858
$ws.Range("J100").Value = 1120.7778
$ws.Range("L100").Value = 2241.5556
$ws.Range("N100").Value = -3323.5556
$ws.Range("H107").Value = 634.6875
$ws.Range("I107").Value = 309.3846
$ws.Range("J107").Value = 857.2632
$ws.Range("K107").Value = 928.1537999999999
$ws.Range("L107").Value = 2571.7896
$ws.Range("M107").Value = 991.8462000000001
$ws.Range("N107").Value = -6411.7896
$ws.Range("H126").Value = 2567.8518
$ws.Range("I126").Value = 2295.5557
$ws.Range("K126").Value = 6886.6671
$ws.Range("M126").Value = -4416.6671
$ws.Range("H132").Value = 6001.1665
$ws.Range("I132").Value = 7875.5
$ws.Range("K132").Value = 23626.5
$ws.Range("M132").Value = -21096.5

